$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
    }
    return $null
}

# --- Reposition existing shapes ---
# NOTE: Shape.Left / Shape.Top round-trip through a single-precision (float32)
# point value internally (1 pt = 12700 EMU), so a naive "emu / 12700" can land
# one float32 ULP below the intended value and truncate to the wrong EMU when
# written back out. The literals below are the verified float32 point values
# that reproduce the exact target EMU offsets from the authoritative OOXML.

# TextBox 5 (id 6): x 3200400 -> 2971800
$sh = Get-ShapeByName $s "TextBox 5"
$sh.Left = 234.0

# TextBox 26 (id 27): x 3148014 -> 2919414
$sh = Get-ShapeByName $s "TextBox 26"
$sh.Left = 229.8751220703125

# Rectangle 14 (id 15): x 3200400 -> 2971800
$sh = Get-ShapeByName $s "Rectangle 14"
$sh.Left = 234.0

# Rectangle 28 (id 29): x 3193257 -> 2964657
$sh = Get-ShapeByName $s "Rectangle 28"
$sh.Left = 233.43756103515625

# Rectangle 15 (id 16): x 5486400 -> 5257800
$sh = Get-ShapeByName $s "Rectangle 15"
$sh.Left = 414.0

# --- Delete the three decorative Chevron shapes ---
$sh = Get-ShapeByName $s "Chevron 33"
if ($sh) { $sh.Delete() }

$sh = Get-ShapeByName $s "Chevron 34"
if ($sh) { $sh.Delete() }

$sh = Get-ShapeByName $s "Chevron 35"
if ($sh) { $sh.Delete() }

# --- Move the three Plus shapes to the right side of the slide ---

# Plus 20 (id 21): x 2650332 -> 8153400, y 1752601 -> 1778290
$sh = Get-ShapeByName $s "Plus 20"
$sh.Left = 642.0
$sh.Top = 140.02284240722656

# Plus 37 (id 38): x 2650332 -> 8153390
$sh = Get-ShapeByName $s "Plus 37"
$sh.Left = 641.999267578125

# Plus 38 (id 39): x 2650332 -> 8153380, y 4269647 -> 4198347
# Also recolor its gradient fill from bg1-based grays to an FFC000 (orange) gradient.
$sh = Get-ShapeByName $s "Plus 38"
$sh.Left = 641.9984741210938
$sh.Top = 330.5785217285156

$fill = $sh.Fill
$grad = $fill.GradientStops
# 49407 (0xC0FF -> BBGGRR) = RGB(255,192,0) = FFC000
$grad.Item(1).Color.RGB = 49407
$grad.Item(1).Position = 0
$grad.Item(2).Color.RGB = 49407
$grad.Item(2).Position = 0.57
$grad.Item(3).Color.RGB = 49407
$grad.Item(3).Position = 1
$fill.RotateWithObject = 0

# --- Delete the other decorative / helper shapes that were dropped from the slide ---
$sh = Get-ShapeByName $s "Rounded Rectangle 23"
if ($sh) { $sh.Delete() }

$sh = Get-ShapeByName $s "Multiply 25"
if ($sh) { $sh.Delete() }

$sh = Get-ShapeByName $s "Multiply 27"
if ($sh) { $sh.Delete() }

$sh = Get-ShapeByName $s "Multiply 29"
if ($sh) { $sh.Delete() }

$sh = Get-ShapeByName $s "Rectangle 36"
if ($sh) { $sh.Delete() }
